$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns H:K
$ws.Range("H1").Value = "tanggal_spk"
$ws.Range("I1").Value = "group"
$ws.Range("J1").Value = "progress"
$ws.Range("K1").Value = "area"

# New data row 2
$ws.Range("H2").Value = "08 Januari 2025"
$ws.Range("I2").Value = "Z"
$ws.Range("J2").Value = 0.5
$ws.Range("J2").NumberFormat = "0%"
$ws.Range("K2").Value = "MA 3"

# Update existing value: E2 "MA 2" -> "UYE"
$ws.Range("E2").Value = "UYE"

# New data row 3
$ws.Range("H3").Value = "12 Desember 2025"
$ws.Range("I3").Value = "A"
$ws.Range("J3").Value = 0.86
$ws.Range("J3").NumberFormat = "0%"
$ws.Range("K3").Value = "MA 4"

# Match column width seen in target (bestFit-like autofit)
$ws.Columns.Item(8).AutoFit() | Out-Null

# Adjust selection to match recorded end-state
$ws.Range("K5").Select()
